$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 35
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 33
$ws.Range("E2").Value = 0

# Row 3
$ws.Range("B3").Value = 21
$ws.Range("C3").Value = 21
$ws.Range("E3").Value = 0

# Row 4
$ws.Range("C4").Value = 19
$ws.Range("D4").Value = 7
